# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-29, column E) is currently sorted
# descending (2402 .. 2301). Re-sort it ascending (2301 .. 2402), carrying
# each row's "Valor Mora" (column F) value along with its period so the
# odd value (33600, originally on the 2402 row) ends up on whichever row
# holds period 2402 once the table is sorted ascending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow = 29

$pairs = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $period = $ws.Cells.Item($r, 5).Value2
    $val = $ws.Cells.Item($r, 6).Value2
    $pairs += [PSCustomObject]@{ Period = $period; Val = $val }
}

$sorted = $pairs | Sort-Object { [int]$_.Period }

$i = 0
foreach ($p in $sorted) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 5).Value2 = $p.Period
    $ws.Cells.Item($r, 6).Value2 = $p.Val
    $i = $i + 1
}
